$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Footer "date" placeholder: 6/2/2021 -> 7/8/2021
#    (Slide Master + every Custom Layout that carries a Date Placeholder.)
# ---------------------------------------------------------------------------
$newDate = "7/8/2021"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 2 (the "figures" slide): reshuffle the "Crop" callout box and the
#    three connector lines that are glued to it, and rename the hemisphere
#    label.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(2)

# Straight Connector 135 (glued to the Crop textbox, idx=3)
$conn1 = $slide.Shapes.Item(5)
$conn1.Left = 488.112060546875
$conn1.Width = 241.10733032226562

# TextBox 102 ("Crop" / hemisphere caption) - move/widen, switch to wrapping
# text, and update the caption itself.
$tb = $slide.Shapes.Item(23)
$origHeight = $tb.Height
$tb.Left = 350.12750244140625
$tb.Width = 137.9845733642578
$tb.TextFrame.WordWrap = -1

$tr = $tb.TextFrame.TextRange
$len = $tr.Length
$caption = $tr.Characters(6, $len - 5)
$caption.Text = "West Northern-Hemisphere"

# Re-wrapping the longer caption at the new width nudges the autofit shape's
# height; the authored change kept the box's height as-is, so restore it.
$tb.Height = $origHeight

# Straight Connector 103 (glued to the Crop textbox, idx=2, flipped)
$conn2 = $slide.Shapes.Item(24)
$conn2.Left = 419.0309753417969
$conn2.Width = 0.08881890028715134

# Straight Connector 133 (glued to the Crop textbox, idx=1)
$conn3 = $slide.Shapes.Item(40)
$conn3.Width = 141.6968536376953
